# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values, repulled from source
$updates = @{
    2  = 1
    3  = -3
    4  = -1
    6  = 2
    7  = -2
    9  = -2
    10 = -3
    11 = -1
    13 = -5
    16 = 0
    17 = -3
    21 = 1
    23 = -3
    26 = 0
    30 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
